# Ajout du dossier IXISPI
# Updates the "Code plaque" reference in B1 and appends the new IMB rows
# (IXISPI dossier) to the "Liste A3" worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of data (Code IMB, Code plaque, Action, ID_ZN, ID_ZE, ID_RGT)
$data = @(
    @("IMB/26198/C/01RA", "CEMRJ1CO_1013", "Qualif Négo", "SO", "SO", "SO"),
    @("IMB/26198/C/01UB", "CEMRJ1CO_1013", "Qualif Négo", "SO", "SO", "SO"),
    @("IMB/26198/C/01UC", "CEMRJ1CO_1013", "Qualif Négo", "SO", "SO", "SO"),
    @("IMB/26198/C/01UD", "CEMRJ1CO_1013", "Qualif Négo", "SO", "SO", "SO"),
    @("IMB/26198/C/02HU", "CEMRJ1CO_1013", "Qualif Négo", "SO", "SO", "SO"),
    @("IMB/26198/C/02MG", "CEMRJ1CO_1013", "Qualif Négo", "SO", "SO", "SO")
)

# Write the new rows starting at row 4, columns A..F (G stays blank)
$r = 4
foreach ($row in $data) {
    $c = 1
    foreach ($val in $row) {
        $ws.Cells.Item($r, $c).Value2 = $val
        $c = $c + 1
    }
    $r = $r + 1
}

# Update the plaque code shown in the header (B1)
$ws.Range("B1").Value2 = "CEMRJ1CO_1013"

# Apply the same bordered style used elsewhere in the table to the new rows
$ws.Range("A4:G9").Borders.LineStyle = 1
